$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '67.428.70'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.87%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.613.24'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.12%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '590.65'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '150.05'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.44%  '
$ws.Range('E7').Value = '  +0.00%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.552'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.75%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.612.18'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.39%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.124'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.19%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  -1.29%  '
$ws.Range('E13').Value = '  -3.01%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.33'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.66%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.085.46'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000179'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -4.48%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '67.312.99'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.65%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.611.63'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.46%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '367.73'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.62%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '11.03'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('E21').Value = '  -4.07%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('E25').Value = '  -0.11%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.91'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E29').Value = '  -0.15%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '579.36'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.68%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0₃0989'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -5.09%  '
$ws.Range('E32').Value = '  -5.12%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '7.65'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('E34').Value = '  -3.21%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.125'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -6.25%  '
$ws.Range('E37').Value = '  -2.45%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '155.46'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.13%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '19.01'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('E42').Value = '  -2.41%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '16.80'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('E45').Value = '  -0.09%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '154.05'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.38%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0₆0284'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.69'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0784'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('E51').Value = '  +1.95%  '
